# Tutorial 6 attendance sheet update:
#  - Reformat the date strings in column A (rows 3-21) from DD/MM/YYYY to DD-MM-YYYY
#  - Update a handful of attendance tally cells (D/E/G/H) for specific rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: swap "/" for "-" in the date text for rows 3 through 21 ---
# Force text format first so Excel doesn't reinterpret the dashed string
# (e.g. "01-08-2022") as an actual date serial, then restore the default
# "Normal" style so no stray number-format sticks to the cell.
for ($row = 3; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $newDate = $cell.Value2.Replace("/", "-")
    $cell.NumberFormat = "@"
    $cell.Value = $newDate
    $cell.Style = "Normal"
}

# --- Attendance tally updates ---
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D7").Value = 1
$ws.Range("G7").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0
